$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sponza" (sheet2.xml) - add a new column K ("v1419") with performance
# figures for the new build, mirroring the existing H/I/J column pattern.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sponza")

$ws2.Range("K1").Value = "v1419"

$ws2.Range("K2").Value = 8415
$ws2.Range("K3").Value = 8327
$ws2.Range("K4").Value = 8288
$ws2.Range("K5").Value = 8295
$ws2.Range("K6").Value = 8285
$ws2.Range("K7").Value = 8292
$ws2.Range("K8").Value = 8333
$ws2.Range("K9").Value = 8308
$ws2.Range("K10").Value = 8271
$ws2.Range("K11").Value = 8272

$ws2.Range("K12").Formula = "=AVERAGE(K2:K11)"
$ws2.Range("K13").Formula = "=_xlfn.VAR.S(K2:K11)"
$ws2.Range("K14").Formula = "=1-_xlfn.T.TEST(J2:J11,K2:K11,2,3)"
$ws2.Range("K15").Formula = "=J12/K12"
$ws2.Range("K16").Formula = "=B12/K12"

# Copy the formatting (styles) from column J onto the new column K, so the
# header/body/summary rows keep the same look (centered header, numeric
# body, wrapped summary rows, ...).
$ws2.Range("J1:J16").Copy() | Out-Null
$ws2.Range("K1:K16").PasteSpecial(-4122) | Out-Null

# Extend the conditional formatting that highlights the DIFF ACCEPT /
# speed-up rows so it also covers the new column K.
$fcs2 = $ws2.Range("B15:J16").FormatConditions
for ($i = 1; $i -le $fcs2.Count; $i++) {
    $fcs2.Item($i).ModifyAppliesToRange($ws2.Range("B15:K16")) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "ComplexMesh" (sheet3.xml) - add the matching new column J ("v1419").
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ComplexMesh")

$ws3.Range("J1").Value = "v1419"

$ws3.Range("J2").Value = 6422
$ws3.Range("J3").Value = 6327
$ws3.Range("J4").Value = 6296
$ws3.Range("J5").Value = 6329
$ws3.Range("J6").Value = 6281
$ws3.Range("J7").Value = 6304
$ws3.Range("J8").Value = 6366
$ws3.Range("J9").Value = 6346
$ws3.Range("J10").Value = 6360
$ws3.Range("J11").Value = 6326
$ws3.Range("J12").Value = 6277

$ws3.Range("J13").Formula = "=_xlfn.VAR.S(J2:J11)"
$ws3.Range("J14").Formula = "=1-_xlfn.T.TEST(I2:I11,J2:J11,2,3)"
$ws3.Range("J15").Formula = "=I12/J12"
$ws3.Range("J16").Formula = "=B12/J12"

# Copy formatting from column I onto the new column J.
$ws3.Range("I1:I16").Copy() | Out-Null
$ws3.Range("J1:J16").PasteSpecial(-4122) | Out-Null

# Extend conditional formatting to cover the new column J.
$fcs3 = $ws3.Range("B15:I16").FormatConditions
for ($i = 1; $i -le $fcs3.Count; $i++) {
    $fcs3.Item($i).ModifyAppliesToRange($ws3.Range("B15:J16")) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet activation / selection - the commit moves the active tab from
# "Sponza" to "ComplexMesh" and updates each sheet's last-selected cell.
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("I22").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("J20").Select() | Out-Null
